$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the three new rows (bootstrapping results) --------------------
# Existing layout (before):
#   row2 BK -> PK
#   row3 D  -> PK
#   row4 P  -> PK
#   row5 WB -> PK
# New layout (after):
#   row2 BK  -> PK   (unchanged)
#   row3 BKD -> PK   (NEW)
#   row4 D   -> PK   (was row3)
#   row5 P   -> PK   (was row4)
#   row6 PD  -> PK   (NEW)
#   row7 WB  -> PK   (was row5)
#   row8 WBD -> PK   (NEW)

# Insert row for "BKD -> PK" right after BK -> PK (old row 2), pushing
# everything below it down by one. The inherited style from the row above
# matches the "s=3" data style used by the whole table.
$ws.Rows(3).Insert()
$ws.Range("A3").Value = "BKD → PK"
$ws.Range("B3").Value = 0.036
$ws.Range("C3").Value = 0.036
$ws.Range("D3").Value = 0.993
$ws.Range("E3").Value = 0.161
$ws.Range("F3").Value = "Tidak signifikan"

# Insert row for "PD -> PK" after "P -> PK" (now row 5 post-insert above).
$ws.Rows(6).Insert()
$ws.Range("A6").Value = "PD → PK"
$ws.Range("B6").Value = -0.044
$ws.Range("C6").Value = 0.041
$ws.Range("D6").Value = 1081
$ws.Range("E6").Value = 0.14
$ws.Range("F6").Value = "Tidak signifikan"

# Append row for "WBD -> PK" after "WB -> PK" (now row 7 post-insert above).
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "WBD → PK"
$ws.Range("B8").Value = 0.008
$ws.Range("C8").Value = 0.031
$ws.Range("D8").Value = 0.263
$ws.Range("E8").Value = 0.397
$ws.Range("F8").Value = "Tidak signifikan"

# --- Add a thin border around all data rows (matches borderId=1) ----------
$ws.Range("A2:F8").Borders.LineStyle = 1

# --- Auto column widths (explicit widths captured on save) ----------------
$ws.Columns(1).ColumnWidth = 10.833333333333334
$ws.Columns(2).ColumnWidth = 24.0
$ws.Columns(3).ColumnWidth = 7.166666666666667
$ws.Columns(4).ColumnWidth = 14.333333333333334
$ws.Columns(5).ColumnWidth = 9.5
$ws.Columns(6).ColumnWidth = 20.333333333333332
